# Add the season record (Wins / Losses / Ties) columns to the roster sheet.
# The old scraper only pulled team statistics, not the win/loss/tie record,
# so this adds three new columns (AD, AE, AF) with the team's season record
# repeated on every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same as the rest of row 1 (bold, centered,
# bordered) by copying the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record, same for every player on the roster.
$wins = 95
$losses = 68
$ties = 0

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 52) { $lastRow = 52 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
